$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("CDCF-PMpPDOU")
$ws3 = $wb.Worksheets.Item("CDCF-FTMpFDOU")

# --- About sheet updates ---

# Row 1-2 header text changed (moved from the old "EU model units" block)
$wsAbout.Range("A1").Value = "CDCF Passenger Miles per Passenger Distance Output Unit"
$wsAbout.Range("A2").Value = "CDCF Freight Ton Miles per Freight Distance Output Unit"

# Row 8 no longer uses the secondary bold-ish style (back to default/Normal)
$wsAbout.Range("A8").Style = "Normal"

# Rows 11-13: replace US wording with EU wording
$wsAbout.Range("A11").Value = "For the EU. model, the desired output units are:"
$wsAbout.Range("A12").Value = "trillion passenger-km"
$wsAbout.Range("A13").Value = "trillion freight ton-km"

# New rows 15-16: conversion factor section with shaded header row
$rngHeader = $wsAbout.Range("A15:B15")
$rngHeader.Interior.ThemeColor = 2
$rngHeader.Interior.TintAndShade = -0.249977111117893
$wsAbout.Range("A15").Value = "Relevant Conversion Factors"
$wsAbout.Range("A15").Font.Bold = $true

$wsAbout.Range("A16").Value = "miles to km"
$wsAbout.Range("B16").Value = 1.60934

# --- CDCF-PMpPDOU / CDCF-FTMpFDOU: multiply by the new miles->km factor ---
$ws2.Range("B2").Formula = "=10^12*About!B16"
$ws3.Range("B2").Formula = "=10^12*About!B16"

# --- Active sheet / selection bookkeeping to mirror the saved view state ---
[void]$wsAbout.Activate()
[void]$wsAbout.Range("B17").Select()

[void]$ws2.Activate()
[void]$ws2.Range("B2").Select()

[void]$ws3.Activate()
[void]$ws3.Range("H13").Select()

Write-Host "done"
